# Update the "Förändrad" date column (C2:C5) from serial 45243 (2023-11-13)
# to serial 45244 (2023-11-14) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
